$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for the "R10" rule row (cell E8) from
# "Good Morning" to "GIT UPDATE".
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the edited cell as the active selection, matching the saved
# workbook view state.
$ws.Range("E8").Select()
